$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F4").Value = 1794
$ws.Range("F5").Value = 3348
$ws.Range("F6").Value = 1144
$ws.Range("F7").Value = 2258
$ws.Range("F8").Value = 2182
$ws.Range("F9").Value = 1144
$ws.Range("F10").Value = 619
$ws.Range("F12").Value = 1705
$ws.Range("F13").Value = 418
$ws.Range("F15").Value = 49
$ws.Range("F16").Value = 312
$ws.Range("F17").Value = 254
$ws.Range("F18").Value = 1614
$ws.Range("F19").Value = 275
$ws.Range("F20").Value = 1333
$ws.Range("F21").Value = 755
$ws.Range("F22").Value = 280
$ws.Range("F23").Value = 635
$ws.Range("F24").Value = 12423
$ws.Range("F25").Value = 12463
$ws.Range("F27").Value = 716
$ws.Range("F28").Value = 6
$ws.Range("F29").Value = 263
$ws.Range("F36").Value = 629

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F6").Value = 1
$ws.Range("F9").Value = 44
$ws.Range("F10").Value = 54

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F3").Value = 113

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F5").Value = 1794
$ws.Range("F6").Value = 3348
$ws.Range("F7").Value = 1144
$ws.Range("F8").Value = 2258
$ws.Range("F9").Value = 2182
$ws.Range("F10").Value = 1144
$ws.Range("F11").Value = 619
$ws.Range("F12").Value = 113
$ws.Range("F14").Value = 1705
$ws.Range("F15").Value = 418
$ws.Range("F18").Value = 49
$ws.Range("F20").Value = 312
$ws.Range("F22").Value = 254
$ws.Range("F23").Value = 1614
$ws.Range("F24").Value = 275
$ws.Range("F25").Value = 1333
$ws.Range("F26").Value = 755
$ws.Range("F27").Value = 280
$ws.Range("F29").Value = 635
$ws.Range("F30").Value = 12423
$ws.Range("F31").Value = 12463
$ws.Range("F33").Value = 716
$ws.Range("F34").Value = 6
$ws.Range("F35").Value = 263
$ws.Range("F38").Value = 1
$ws.Range("F44").Value = 44
$ws.Range("F46").Value = 629
$ws.Range("F47").Value = 54
